$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2 "Latest Handback DateTime"
# both move from 2016-08-18 15:04:15 -> 2016-08-18 15:05:06
$wsOverview.Range("G2").Value = "2016-08-18 15:05:06"
$wsDeDe.Range("H2").Value     = "2016-08-18 15:05:06"

# zh-cn!H2 "Latest Handoff Datetime" moves from 2016-08-18 15:04:00 -> 2016-08-18 15:04:55
$wsZhCn.Range("H2").Value = "2016-08-18 15:04:55"

# --- Column width updates: the "Status" columns shrink from ~29.98 to ~17.22 ---
$wsOverview.Range("E1").ColumnWidth = 16.333333
$wsOverview.Range("F1").ColumnWidth = 16.333333
$wsZhCn.Range("C1").ColumnWidth = 16.333333
$wsDeDe.Range("C1").ColumnWidth = 16.333333
